# Regenerate the "K" (strikeouts) column (column G) values in the save_data
# sheet. The source stats were re-pulled (K used instead of Strike#), so the
# resulting strikeout counts differ from what was previously stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 6
    3  = 6
    4  = 7
    5  = 4
    6  = 8
    7  = 6
    8  = 12
    9  = 7
    10 = 4
    11 = 5
    12 = 8
    13 = 10
    14 = 1
    15 = 4
    16 = 7
    17 = 5
    18 = 5
    19 = 7
    20 = 7
    21 = 1
    22 = 8
    23 = 8
    24 = 8
    25 = 13
    26 = 8
    27 = 9
    28 = 10
    29 = 5
    30 = 3
    31 = 7
    32 = 7
    33 = 8
    34 = 7
    35 = 4
    36 = 2
    37 = 4
    38 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
